# edit.ps1 -- apply the three text edits described by the commit diff.
#
#  1. Slide 11, content placeholder, 2nd paragraph:
#       "Verification is \u2026"
#       -> "Verification is through nodes in the network (we will talk about it in a moment)"
#
#  2. Slide 9, title:
#       "How are we identified with this coin?"
#       -> "How are we identified" + "(?)" (yellow-highlighted) + " with this coin?"
#
#  3. Slide 9, content placeholder, 1st paragraph:
#       "We are identified by our public key"
#       -> "We are " + "identified" (yellow-highlighted) + " by our public key"

$p = $ppt.ActivePresentation

# Yellow highlight, same RGB() value PowerPoint uses for FFFF00 (R=255,G=255,B=0).
$Yellow = 65535

# --- 1. Slide 11: expand the "Verification is ..." bullet -------------------
$slide11 = $p.Slides.Item(11)
$body11 = $slide11.Shapes.Item(2).TextFrame.TextRange
$para2 = $body11.Paragraphs(2, 1)
$body11.Characters($para2.Start, $para2.Length).Text = "Verification is through nodes in the network (we will talk about it in a moment)"

# --- 2. Slide 9 title: "How are we identified with this coin?" --------------
$slide9 = $p.Slides.Item(9)
$title9 = $slide9.Shapes.Item(1).TextFrame.TextRange

# Split after "How are we identified" (21 characters) by inserting the new
# "(?)" run right after it; the trailing " with this coin?" is left as-is.
$lead = $title9.Characters(1, 21)
[void]$lead.InsertAfter("(?)")

$qmark = $title9.Characters(22, 3)
$qmark.Font.Highlight.RGB = $Yellow

# --- 3. Slide 9 content placeholder: "We are identified by our public key" --
$body9 = $slide9.Shapes.Item(2).TextFrame.TextRange
$para1 = $body9.Paragraphs(1, 1)
$word = $para1.Characters(8, 10)
$word.Font.Highlight.RGB = $Yellow
